# Update cryptos list values (prices, % volumes, and two swapped rows) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (some look numeric, e.g. "314.80" or
# "0.0000154"); force those cells to Text format first so Excel does not silently
# coerce them into numbers and drop significant trailing zeros / reformat them.
$ws.Range("D2").Value = "62.667.06"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "3.441.07"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.75"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.09"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.694"
$ws.Range("E9").Value = "  +3.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  +11.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.05"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.45"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.86"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "3.441.30"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").Value = "62.669.50"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.57"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000154"
$ws.Range("E19").Value = "  +17.48%  "
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "314.80"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.84"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.18"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.76"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.76"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.16"
$ws.Range("E27").Value = "  -3.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.83"
$ws.Range("E28").Value = "  +5.09%  "
$ws.Range("E29").Value = "  +7.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "44.57"
$ws.Range("E30").Value = "  +7.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.174"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.115"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.39"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0484"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.76"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.325"
$ws.Range("E38").Value = "  +15.04%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.97"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "142.60"
$ws.Range("E41").Value = "  +4.74%  "
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.87"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.35"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("D48").Value = "2.113.56"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.97"
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.10"
$ws.Range("E51").Value = "  +29.94%  "
